$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 646, shifting rows 646:687 down to 647:688.
$ws.Rows.Item(646).Insert()

# Fill in the newly inserted row with the new record:
#   2026/01/18  日  4  19
# Force column A to be stored as text (matching the rest of the date
# column, which holds text like "2026/01/18" rather than a real date
# serial number), then restore the "Normal" style so no stray number
# format is left behind on the cell.
$ws.Range("A646").NumberFormat = "@"
$ws.Range("A646").Value = "2026/01/18"
$ws.Range("A646").Style = "Normal"
$ws.Range("B646").Value = "日"
$ws.Range("C646").Value = 4
$ws.Range("D646").Value = 19
